# Apply the edits described by the commit: renumber some Sim numbers,
# append new rows of Vina / Tu quy sims to Sheet1, append git workflow
# notes to Sheet2, and update the active sheet / selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------
# Sheet1: renumber rows 13-24 (column A) to the new phone numbers
# ---------------------------------------------------------------
$renumber = @{
    13 = 987580501
    14 = 987580502
    15 = 987580503
    16 = 987580504
    17 = 987580505
    18 = 987580506
    19 = 987580507
    20 = 987580508
    21 = 987580509
    22 = 987580510
    23 = 987580511
    24 = 987580512
}
foreach ($r in $renumber.Keys) {
    $ws1.Range("A" + $r).Value = $renumber[$r]
}

# ---------------------------------------------------------------
# Sheet1: append new rows 32-49 with Vina sims
# (seed the "Vina" / "Tứ quý" shared strings first so they land at
# shared-string indices 17 and 18, matching the target workbook)
# ---------------------------------------------------------------
$ws1.Range("C32").Value = "Vina"
$ws1.Range("D39").Value = "Tứ quý"

$newRows = @(
    @(32, 987585574, "Lộc phát"),
    @(33, 987585575, "Lộc phát"),
    @(34, 987585576, "Lộc phát"),
    @(35, 987585577, "Lộc phát"),
    @(36, 987585578, "Lộc phát"),
    @(37, 987585579, "Lộc phát"),
    @(38, 987585580, "Lộc phát"),
    @(39, 987585581, "Tứ quý"),
    @(40, 987585582, "Tứ quý"),
    @(41, 987585583, "Tứ quý"),
    @(42, 987585584, "Tứ quý"),
    @(43, 987585585, "Tứ quý"),
    @(44, 987585586, "Tứ quý"),
    @(45, 987585587, "Tứ quý"),
    @(46, 987585588, "Tứ quý"),
    @(47, 987585589, "Tứ quý"),
    @(48, 987585590, "Tứ quý"),
    @(49, 987585591, "Tứ quý")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws1.Range("A" + $r).Value = $row[1]
    $ws1.Range("B" + $r).Value = 1200000
    $ws1.Range("C" + $r).Value = "Vina"
    $ws1.Range("D" + $r).Value = $row[2]
    $ws1.Range("E" + $r).Value = "dễ nhớ"
}

# ---------------------------------------------------------------
# Sheet2: append the git workflow notes in rows 4-6
# ---------------------------------------------------------------
$ws2.Range("A4").Value = "git add ."
$ws2.Range("A5").Value = 'git commit -m "Chỉnh sửa [mô tả phần chỉnh sửa]"'
$ws2.Range("A6").Value = "git push origin main"

# ---------------------------------------------------------------
# Update selection / active sheet state:
#  - Sheet1 keeps a selection on F27 but is no longer the active tab
#  - Sheet2 becomes the active tab with A4:A6 selected
# ---------------------------------------------------------------
$ws1.Range("F27").Select()
$ws2.Select()
$ws2.Range("A4:A6").Select()
